$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report title (shared string content change: "Trends.Earth..." -> "Misland...")
$ws.Range("A1").Value = "Misland potential carbon removals from restoration summary table"

# Clear the old "For more information on Trends.Earth..." contact blurb cell
$ws.Range("A20").Value = ""

# Move the active selection to A20, matching the saved cursor position
[void]$ws.Range("A20").Select()
